$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.276.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.720.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4729"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2625"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06197"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.716.28"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07055"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5978"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.434"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9999"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9997"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.284.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006811"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.95%  "

$ws.Range("E20").Value = "  +0.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.935.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.532"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.712"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.248"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.53%  "

$ws.Range("E26").Value = "  +1.20%  "

$ws.Range("E27").Value = "  +3.31%  "

$ws.Range("E28").Value = "  +0.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.945"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.690"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07796"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04506"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.614"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9791"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6207"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9324"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "114.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +18.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.437"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.928"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9993"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.620"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +15.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.01481"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3827"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1180"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.337"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05261"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.773"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3385"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.68%  "

$ws.Range("E51").Value = "  +2.05%  "

